# Additional scraping: add a "Player Info" sheet, and replace the full
# howstat match-card URL columns with a short numeric MATCH_CODE column
# on both the "ODI Batting" and "ODI Bowling" sheets.

$wb = $excel.ActiveWorkbook

$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

# Match codes, in the same row order as they already appear on both sheets.
$matchCodes = @("4086","4182","4183","4186","4187","4188","4206","4247","4261","4264","4488","4491")

# ---------------------------------------------------------------------
# 1) "ODI Batting": MATCH_CARD_LINK (col D) -> MATCH_CODE, URL -> numeric code
# ---------------------------------------------------------------------
$wsBatting.Range("D1").Value = "MATCH_CODE"
for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $row = $i + 2
    $cell = $wsBatting.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$i]
}

# ---------------------------------------------------------------------
# 2) "ODI Bowling": MATCH_CARD_LINK (col B) -> MATCH_CODE, URL -> numeric code
# ---------------------------------------------------------------------
$wsBowling.Range("B1").Value = "MATCH_CODE"
for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $row = $i + 2
    $cell = $wsBowling.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$i]
}

# ---------------------------------------------------------------------
# 3) Insert a new "Player Info" sheet before "ODI Batting" so the final
#    sheet order is: Player Info, ODI Batting, ODI Bowling
# ---------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Add($wsBatting)
$wsInfo.Name = "Player Info"

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsInfo.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$headerRange = $wsInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$wsInfo.Range("A2").NumberFormat = "@"
$wsInfo.Range("A2").Value = "4707"
$wsInfo.Range("B2").Value = "Pieter Willem Adriaan Mulder"
$wsInfo.Range("C2").Value = "Right Handed"
$wsInfo.Range("D2").Value = "Right Arm Medium"

$wsInfo.Range("A1").Select()
